$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Draft")

$ws.Range("A2").Value = "Aidan || 7135Z"
$ws.Range("A3").Value = "Michael || 676D"
$ws.Range("A4").Value = "Jeremiah || OBSR || 676V"
$ws.Range("A5").Value = "Edward || 929"
$ws.Range("A6").Value = "Donald || 20850V"
$ws.Range("A7").Value = "Charlie || 929K"
$ws.Range("A8").Value = "Jerry || 20850Z"
$ws.Range("A9").Value = "Rowan || 593C"
$ws.Range("A10").Value = "James || 9080S"
$ws.Range("A11").Value = "Jack || 934Z"
$ws.Range("A12").Value = "Leah || 929T"
$ws.Range("A13").Value = "Chad || 53E"
$ws.Range("A14").Value = "Amelia || 20850A"
